# Appends three new "created leaders" log entries (rows 186-188) to the
# log sheet, matching the new rows added by the app ("neue Funktionen,
# ETO Dialog für den Export"). The sheet already has pre-formatted empty
# rows reserved up to row 1000, so we just fill in the next three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 186
$ws.Range("A186").Value = "2025-10-02T11:28:41"
$ws.Range("B186").Value = "Untitled"
$ws.Range("C186").Value = "37aed359-c26a-4879-a3b0-be72d0102e6c"
$ws.Range("D186").Value = "zargentuere"
$ws.Range("E186").Value = "Standard 1:10 Zargenbeschriftung"
$ws.Range("F186").Value = "zargentuere_zimmer_M150.csv"

# Row 187
$ws.Range("A187").Value = "2025-10-02T13:23:55"
$ws.Range("B187").Value = "Haus_A.3dm"
$ws.Range("C187").Value = "902cfd9e-9c4e-43ed-a10e-439db871fb13"
$ws.Range("D187").Value = "schiebetuere"
$ws.Range("E187").Value = "Standard 1:10 Schiebetürbeschriftung"
$ws.Range("F187").Value = "schiebetuere.csv"

# Row 188
$ws.Range("A188").Value = "2025-10-02T15:18:23"
$ws.Range("B188").Value = "Haus_A.3dm"
$ws.Range("C188").Value = "f6f1031a-2245-4d3d-9f7c-9d911d750bbb"
$ws.Range("D188").Value = "schiebetuere"
$ws.Range("E188").Value = "Standard 1:10 Schiebetürbeschriftung"
$ws.Range("F188").Value = "schiebetuere.csv"
